# Fix config file handling
# Appends one new trailing data row (row 83) to each of the four sensor-data
# worksheets, mirroring the existing row layout (A: timestamp text,
# B/C/D/E: hex string fields, F: numeric, G: big-integer-as-text, H/I: numeric).

$wb = $excel.ActiveWorkbook

# New row 83 values per worksheet, in workbook (tab) order:
#   1 ROW35-FE-LIFTER, 2 ROW35-MID-LIFTER, 3 ROW02-FE-LIFTER, 4 ROW02-MID-LIFTER
$newRows = @(
    @{
        A = "2025-03-07 18:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        A = "2025-03-07 18:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        A = "2025-03-07 18:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        A = "2025-03-07 18:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $row = $newRows[$i - 1]
    $r = 83

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F

    # Column G holds a 24-digit number far beyond double precision in the
    # source data; it is stored as text, so force text formatting before
    # assignment so the long digit string is preserved exactly.
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $row.G

    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
